$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.304.36'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.868.45'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4704'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2869'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06580'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08013'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").Value = '1.872.32'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.121'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6845'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '269.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.41%  '
$ws.Range("D17").Value = '30.341.22'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007634'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '2.119.29'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.285'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.215'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.409'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.949'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09926'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.364'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.462'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.065'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04721'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7006'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.710'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.638'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.280'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.81%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8431'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4169'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9993'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.174'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.046'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '907.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05703'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.01%  '

Write-Output "Applied changes"